$d = $word.ActiveDocument

# Locate the "Unsubscribe" hyperlink (added at the end of the document) and
# remove it entirely -- including its display text -- while leaving the
# (now empty) paragraph that contained it intact.
$target = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "Unsubscribe") {
        $target = $h
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Trim the trailing paragraph mark off the range so only the hyperlink's
    # run content (the "Unsubscribe" text) is removed, leaving the paragraph
    # itself (and its indentation formatting) in place but empty.
    $delRange = $d.Range($r.Start, $r.End - 1)
    $delRange.Delete()
}
